# Weekly driver report update for 2025-04-20
# Updates the "Bad Drivers" summary rows (3-6) and refreshes the
# "Good Drivers" table (rows 14-22) with the latest roaming stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a value into a cell without Excel's automatic
#     "looks like a date" coercion kicking in (used for the
#     Driver Vintage column, which stores ISO date strings as plain
#     text, not real dates). ---
function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

# ---------------------------------------------------------------
# "Bad Drivers" block (rows 3-6)
# ---------------------------------------------------------------
$ws.Range("C3").Value = 414
$ws.Range("D3").Value = 95.3

$ws.Range("C4").Value = 730
$ws.Range("D4").Value = 97.59999999999999

$ws.Range("D5").Value = 98.8

$ws.Range("C6").Value = 1192

# ---------------------------------------------------------------
# "Good Drivers" block (rows 14-22): Adapter-Driver, Total Samples,
# Good Roaming Calculation (%), Driver Vintage
# ---------------------------------------------------------------

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B14").Value = 445055
$ws.Range("D14").Value = 99.90000000000001
Set-TextValue $ws.Range("E14") "2024-11-10"

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B15").Value = 10661
$ws.Range("D15").Value = 100
Set-TextValue $ws.Range("E15") "2022-08-29"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B16").Value = 14239
$ws.Range("D16").Value = 100
Set-TextValue $ws.Range("E16") "2022-05-23"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B17").Value = 265400
$ws.Range("D17").Value = 99.90000000000001
Set-TextValue $ws.Range("E17") "2022-05-01"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B18").Value = 77849
$ws.Range("D18").Value = 99.90000000000001
Set-TextValue $ws.Range("E18") "2021-08-18"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B19").Value = 34244
$ws.Range("D19").Value = 100
Set-TextValue $ws.Range("E19") "2021-04-27"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B20").Value = 59673
$ws.Range("D20").Value = 100
Set-TextValue $ws.Range("E20") "2020-08-05"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B21").Value = 113652
$ws.Range("D21").Value = 100
Set-TextValue $ws.Range("E21") "2020-01-06"

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B22").Value = 56018
$ws.Range("D22").Value = 100
Set-TextValue $ws.Range("E22") "2019-12-14"
